# Add 2 more DM levels to show variation in output.
#
# The "Slurry" sheet originally has 4 rows (man.dm = 5.9) for the acid
# doses 0 / 11 / 3.4 / 7.5 kg/t. The authored change:
#   1. Lowers the existing DM level (C2:C5) from 5.9 to 5.1.
#   2. Re-adds the original DM=5.9 block as rows 6-9.
#   3. Adds a brand-new DM=6.9 block as rows 10-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slurry")

# --- Step 1: lower the first DM block from 5.9 to 5.1 ---
$ws.Range("C2").Value = 5.1
$ws.Range("C3").Value = 5.1
$ws.Range("C4").Value = 5.1
$ws.Range("C5").Value = 5.1

# Data describing each of the 4 rows within one DM block.
$acidLabels = @("0 kg/t", "11 kg/t", "3.4 kg/t", "7.5 kg/t")
$phFormulas = @($null, "=7.9-1.38", "=7.9-0.8187", "=7.9-1.11")

function Add-DMBlock([int]$startRow, [double]$dmValue) {
    for ($i = 0; $i -lt 4; $i++) {
        $r = $startRow + $i

        $ws.Cells.Item($r, 1).Value = "Afgasset biomasse"
        $ws.Cells.Item($r, 2).Value = $acidLabels[$i]
        $ws.Cells.Item($r, 3).Value = $dmValue

        if ($i -eq 0) {
            $ws.Cells.Item($r, 4).Value = 7.9
            $ws.Cells.Item($r, 4).NumberFormat = "0.00"
        } else {
            $ws.Cells.Item($r, 4).Formula = $phFormulas[$i]
            $ws.Cells.Item($r, 4).NumberFormat = "0.00"
            # OLE COLORREF (0xBBGGRR) matching the workbook's existing
            # "computed ph" font colour FF3465A4 (RGB 0x34,0x65,0xA4).
            $ws.Cells.Item($r, 4).Font.Color = 10773812
        }
    }
}

# --- Step 2: re-add the original DM=5.9 block as rows 6-9 ---
Add-DMBlock 6 5.9

# --- Step 3: add a new DM=6.9 block as rows 10-13 ---
Add-DMBlock 10 6.9

# Mirror the authored selection change.
$ws.Range("D16").Select()
